# Apply cryptos.xlsx data refresh per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.213.71'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '2.267.20'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.27'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.18'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.83%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.494'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.80%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.20'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.85%  '
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '48.30'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -6.12%  '
$ws.Range('E13').Value = '  +0.83%  '
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.62'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.63%  '
$ws.Range('D16').Value = '2.619.53'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').Value = '2.269.88'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.784'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.05%  '
$ws.Range('D19').Value = '42.131.68'
$ws.Range('E19').Value = '  -0.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.66'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.71%  '
$ws.Range('D21').Value = '0.0₃0889'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.00'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.33'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.46%  '
$ws.Range('E24').Value = '  +0.57%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.46'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.95%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.92'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -3.92%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.21'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.04%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '168.42'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.18'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '33.64'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('E34').Value = '  -2.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.56'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('E36').Value = '  -1.99%  '
$ws.Range('E37').Value = '  -3.87%  '
$ws.Range('E38').Value = '  -3.57%  '
$ws.Range('E39').Value = '  -3.30%  '
$ws.Range('E40').Value = '  -1.46%  '
$ws.Range('E41').Value = '  -1.62%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.72'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.39'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.36%  '
$ws.Range('D44').Value = '1.970.75'
$ws.Range('E44').Value = '  -0.22%  '
$ws.Range('E45').Value = '  -0.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.44'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -6.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.55'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -5.27%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.78'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.59%  '
$ws.Range('D49').Value = '2.492.67'
$ws.Range('E49').Value = '  -0.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '52.29'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -5.22%  '
$ws.Range('E51').Value = '  -0.11%  '
